# Automatic hashcode update - applies new MD5-style hashcodes to the
# rows identified by their code in column A (hashcode value lives in
# column B of the same row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B100" = "85819c9b0ee864700a6fb3abf7b62758"  # 04-040021TM
    "B104" = "afc45b0ea45fcd2114d8102997488408"  # 04-040021TP
    "B113" = "956b266fd844e9f3fca2194ee278fadb"  # 04-040021TC
    "B122" = "d15ca3c8fb72fbbd22db7c2394f28a69"  # 04-040014TC
    "B164" = "0a80cf60deec27272e68c8141fbee685"  # 04-040021A
    "B230" = "a7ccd9496d18261177551264266f67e7"  # 04-040014TP
    "B233" = "380c5e4c6ed05e85df43317f9a0cfa66"  # 04-040014TM
    "B331" = "d9986ed4380897b50d61c0803314de7c"  # 04-040018TP
    "B342" = "052d5b4453144717d9154004c40aed09"  # 04-040018TC
    "B343" = "9c8e173b79f48d63f00af95644862e76"  # 04-040018TM
    "B381" = "ccb51bd55ef71d785c4cbe725d27c184"  # 01-010073A
    "B419" = "930e9bd628ccd09c643cd2b4a4b8cfad"  # 05-0709-070905BTC
    "B458" = "752988414c894035dd2770010236af04"  # 01-010073TP
    "B477" = "d42521fa4802f5f3088dfd72d207e8c7"  # 01-010073TC
    "B619" = "bd09cfb4e9f5a5a1edc58ee2f6cbef23"  # 04-040015TC
    "B623" = "5df9e1ffb7ca51b90d6720532ccfee6f"  # 04-040015TP
    "B628" = "ae8a27b09551a4de674da30e82a0e23c"  # 04-040015TM
    "B779" = "babf3fd530aff2ea45435a4292853ff1"  # 04-040018A
    "B818" = "4c2ed9e49577e877cba8646fab52dc00"  # 04-040015A
    "B831" = "3ebef27ff7385eb5bb0c6c1d9dc07834"  # 04-040014A
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
